# Adds a new "2022-Q3" sheet (with its fund-holding detail table) ahead of
# the existing "2022-Q2" sheet, and updates the "总计" (totals) summary
# sheet with the corresponding new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" summary sheet: shift the existing 7 data rows down by
#    one and insert the new 2022-Q3 totals as the new row 2.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$summaryRows = @(
    @(0, "2022-Q3", 7,  "0.38"),
    @(1, "2022-Q2", 8,  "0.47"),
    @(2, "2022-Q1", 13, "1.46"),
    @(3, "2021-Q4", 17, "4.39"),
    @(4, "2021-Q3", 5,  "0.22"),
    @(5, "2021-Q2", 19, "3.31"),
    @(6, "2021-Q1", 12, "18.69"),
    @(7, "2020-Q4", 1,  "7.73")
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $summaryRows[$i][0]
    $total.Cells.Item($r, 2).Value = $summaryRows[$i][1]
    $total.Cells.Item($r, 3).Value = $summaryRows[$i][2]
    $total.Cells.Item($r, 4).Value = $summaryRows[$i][3]
}

# Row 9 (2020-Q4) is brand new territory on this sheet - copy column A's
# "index" styling (bold + border) down from row 8 so it matches the rest of
# the column.
$total.Cells.Item(8, 1).Copy($total.Cells.Item(9, 1))
$total.Cells.Item(9, 1).Value = $summaryRows[7][0]

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" sheet right before "2022-Q2".
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2)
$newSheet.Name = "2022-Q3"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# NOTE: $q2 becomes a stale reference once a sheet is inserted ahead of it -
# re-resolve it by name before using it again.
$q2 = $wb.Worksheets.Item("2022-Q2")

# Header row (B1:H1) - copy text + styling from the 2022-Q2 sheet so the
# bold/border formatting matches exactly, then fix up the one new column of
# text.
$q2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))

$fundRows = @(
    @(0, "001322", "东吴新趋势价值线灵活配置混合", "1.71", "93.68", "7.76", "0.1327", 6),
    @(1, "580002", "东吴价值成长双动力混合A",       "2.69", "90.84", "3.86", "0.1038", 5),
    @(2, "014376", "东吴新能源汽车股票A",           "0.97", "94.05", "7.64", "0.0741", 6),
    @(3, "001323", "东吴移动互联灵活配置混合A",     "0.61", "93.49", "7.30", "0.0445", 6),
    @(4, "014377", "东吴新能源汽车股票C",           "0.27", "94.05", "7.64", "0.0206", 6),
    @(5, "002170", "东吴移动互联灵活配置混合C",     "0.05", "93.49", "7.30", "0.0036", 6),
    @(6, "011241", "东吴价值成长双动力混合C",       "0.00", "90.84", "3.86", "__NUM_ZERO__", 5)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    # Column A uses the same bold/bordered "index" style as the header and
    # as column A on the other quarterly sheets - copy it across first, then
    # overwrite with this row's value.
    $q2.Cells.Item($r, 1).Copy($newSheet.Cells.Item($r, 1))
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 2).ClearFormats()

    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 3).ClearFormats()

    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 4).ClearFormats()

    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 5).ClearFormats()

    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 6).ClearFormats()

    if ($row[6] -eq "__NUM_ZERO__") {
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $newSheet.Cells.Item($r, 7).NumberFormat = "@"
        $newSheet.Cells.Item($r, 7).Value = $row[6]
        $newSheet.Cells.Item($r, 7).ClearFormats()
    }

    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# Restore the original active sheet/selection (总计, A1) so we don't leave
# the new sheet as the active tab.
[void]$total.Activate()
[void]$total.Range("A1").Select()

